$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New hour labels (column A, rows 2-25) and their updated counts (column B),
# re-sorted in descending order of count as in the updated source data.
$labels = @("08","16","09","10","11","13","17","14","15","24","07","18","20","19","21","22","06","23","12","05","01","03","04","02")
$values = @(232,179,174,172,156,145,135,132,129,128,102,57,56,50,47,42,40,39,27,21,21,19,15,10)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
